$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4120336.2
$ws.Range("I40").Value = 3111.1667
$ws.Range("J40").Value = 7414116.5
$ws.Range("K40").Value = 3111.1667
$ws.Range("L40").Value = 7414116.5
$ws.Range("M40").Value = -2936.1667
$ws.Range("N40").Value = -7414466.5

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H116").Value = 7413.9688
$ws.Range("J116").Value = 6819.5
$ws.Range("L116").Value = 6819.5
$ws.Range("N116").Value = -13703.5

$ws.Range("H135").Value = 75000340
$ws.Range("I135").Value = 26316144
$ws.Range("K135").Value = 236845296
$ws.Range("M135").Value = -236842761

$ws.Range("H137").Value = 2085.9167
$ws.Range("I137").Value = 1463.9474
$ws.Range("J137").Value = 4449.4
$ws.Range("K137").Value = 4391.8422
$ws.Range("L137").Value = 13348.2
$ws.Range("M137").Value = -1841.8422
$ws.Range("N137").Value = -18448.2

$ws.Range("H138").Value = 3925.4707
$ws.Range("I138").Value = 2523
$ws.Range("K138").Value = 7569
$ws.Range("M138").Value = -2429

$ws.Range("H141").Value = 2848.25
$ws.Range("I141").Value = 2848.25
$ws.Range("K141").Value = 8544.75
$ws.Range("M141").Value = -3364.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3382.6086
$ws.Range("I32").Value = 2330.2327
$ws.Range("K32").Value = 2330.2327
$ws.Range("M32").Value = -2043.2327

$ws.Range("H63").Value = 1375.75
$ws.Range("I63").Value = 1001
$ws.Range("J63").Value = 1750.5
$ws.Range("K63").Value = 1001
$ws.Range("L63").Value = 1750.5
$ws.Range("M63").Value = -315
$ws.Range("N63").Value = -3122.5

$ws.Range("H66").Value = 1375.75
$ws.Range("I66").Value = 1001
$ws.Range("J66").Value = 1750.5
$ws.Range("K66").Value = 5005
$ws.Range("L66").Value = 8752.5
$ws.Range("M66").Value = -1573
$ws.Range("N66").Value = -15616.5

$ws.Range("H74").Value = 58827020
$ws.Range("I74").Value = 66669290
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 66669290
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = -66668416
$ws.Range("N74").Value = -11748

$ws.Range("H77").Value = 58827020
$ws.Range("I77").Value = 66669290
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 333346450
$ws.Range("L77").Value = 50000
$ws.Range("M77").Value = -333342082
$ws.Range("N77").Value = -58736

$ws.Range("H103").Value = 65000
$ws.Range("J103").Value = 65000
$ws.Range("L103").Value = 65000
$ws.Range("N103").Value = -67344

$ws.Range("H105").Value = 99000
$ws.Range("J105").Value = 99000
$ws.Range("L105").Value = 99000
$ws.Range("N105").Value = -105988

$ws.Range("H122").Value = 5666
$ws.Range("I122").Value = 4763.3076
$ws.Range("K122").Value = 14289.9228
$ws.Range("M122").Value = -11839.9228

$ws.Range("H132").Value = 4169121.5
$ws.Range("I132").Value = 4547314.5
$ws.Range("K132").Value = 13641943.5
$ws.Range("M132").Value = -13639413.5

$ws.Range("H140").Value = 64995
$ws.Range("J140").Value = 64995
$ws.Range("L140").Value = 64995
$ws.Range("N140").Value = -75355

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1981.8
$ws.Range("I7").Value = 2436
$ws.Range("J7").Value = 1073.4
$ws.Range("K7").Value = 2436
$ws.Range("L7").Value = 1073.4
$ws.Range("M7").Value = -2323
$ws.Range("N7").Value = -1299.4

$ws.Range("H31").Value = 10870.667
$ws.Range("I31").Value = 10870.667
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 10870.667
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -10575.667
$ws.Range("N31").ClearContents()

$ws.Range("H34").Value = 10870.667
$ws.Range("I34").Value = 10870.667
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 10870.667
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -10668.667
$ws.Range("N34").ClearContents()

$ws.Range("H58").Value = 16134265
$ws.Range("I58").Value = 31257518
$ws.Range("J58").Value = 2795.2
$ws.Range("K58").Value = 31257518
$ws.Range("L58").Value = 2795.2
$ws.Range("M58").Value = -31257315
$ws.Range("N58").Value = -3201.2

$ws.Range("H80").Value = 60999.5
$ws.Range("J80").Value = 60999.5
$ws.Range("L80").Value = 60999.5
$ws.Range("N80").Value = -63245.5

$ws.Range("H83").Value = 60999.5
$ws.Range("J83").Value = 60999.5
$ws.Range("L83").Value = 182998.5
$ws.Range("N83").Value = -194230.5

$ws.Range("H132").Value = 166668020
$ws.Range("I132").Value = 200001250
$ws.Range("K132").Value = 600003750
$ws.Range("M132").Value = -600001220

$ws.Range("H134").Value = 35859730
$ws.Range("I134").Value = 50201020
$ws.Range("K134").Value = 150603060
$ws.Range("M134").Value = -150600525

$ws.Range("H136").Value = 16134265
$ws.Range("I136").Value = 31257518
$ws.Range("J136").Value = 2795.2
$ws.Range("K136").Value = 93772554
$ws.Range("L136").Value = 8385.599999999999
$ws.Range("M136").Value = -93770004
$ws.Range("N136").Value = -13485.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 65.09090999999999
$ws.Range("I2").Value = 83.333336
$ws.Range("J2").Value = 58.25
$ws.Range("K2").Value = 500.000016
$ws.Range("L2").Value = 349.5
$ws.Range("M2").Value = -387.000016
$ws.Range("N2").Value = -575.5

$ws.Range("H32").Value = 100366664
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H60").Value = 7049.8
$ws.Range("I60").Value = 166
$ws.Range("K60").Value = 498
$ws.Range("M60").Value = -247

$ws.Range("H131").Value = 2398
$ws.Range("I131").Value = 2068
$ws.Range("J131").Value = 2590.5
$ws.Range("K131").Value = 6204
$ws.Range("L131").Value = 7771.5
$ws.Range("M131").Value = -1164
$ws.Range("N131").Value = -17851.5

$ws.Range("H132").Value = 3678
$ws.Range("I132").Value = 2999.25
$ws.Range("J132").Value = 4583
$ws.Range("K132").Value = 26993.25
$ws.Range("L132").Value = 41247
$ws.Range("M132").Value = -24463.25
$ws.Range("N132").Value = -46307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1628.8
$ws.Range("I97").Value = 1628.8
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1628.8
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1132.8
$ws.Range("N97").ClearContents()

$ws.Range("H105").Value = 59798.25
$ws.Range("J105").Value = 59798.25
$ws.Range("L105").Value = 59798.25
$ws.Range("N105").Value = -66786.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2754.348
$ws.Range("I7").Value = 2793.3809
$ws.Range("K7").Value = 2793.3809
$ws.Range("M7").Value = -2681.3809

$ws.Range("H22").Value = 1913.0625
$ws.Range("I22").Value = 1884.9166
$ws.Range("K22").Value = 1884.9166
$ws.Range("M22").Value = -1589.9166

$ws.Range("H27").Value = 1913.0625
$ws.Range("I27").Value = 1884.9166
$ws.Range("K27").Value = 1884.9166
$ws.Range("M27").Value = -1777.9166

$ws.Range("H46").Value = 1262
$ws.Range("I46").Value = 1470.3334
$ws.Range("J46").Value = 637
$ws.Range("K46").Value = 1470.3334
$ws.Range("L46").Value = 637
$ws.Range("M46").Value = -1282.3334
$ws.Range("N46").Value = -1013

$ws.Range("H100").Value = 11744244
$ws.Range("I100").Value = 12478016
$ws.Range("J100").Value = 3899
$ws.Range("K100").Value = 12478016
$ws.Range("L100").Value = 3899
$ws.Range("M100").Value = -12477475
$ws.Range("N100").Value = -4981

$ws.Range("H101").Value = 27999.666
$ws.Range("J101").Value = 27999.666
$ws.Range("L101").Value = 27999.666
$ws.Range("N101").Value = -34489.666

$ws.Range("H126").Value = 2754.348
$ws.Range("I126").Value = 2793.3809
$ws.Range("K126").Value = 8380.1427
$ws.Range("M126").Value = -5910.1427

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 49160.617
$ws.Range("I81").Value = 53809.367
$ws.Range("J81").Value = 4997.5
$ws.Range("K81").Value = 107618.734
$ws.Range("L81").Value = 9995
$ws.Range("M81").Value = -106557.734
$ws.Range("N81").Value = -12117

$ws.Range("H84").Value = 49160.617
$ws.Range("I84").Value = 53809.367
$ws.Range("J84").Value = 4997.5
$ws.Range("K84").Value = 538093.6699999999
$ws.Range("L84").Value = 49975
$ws.Range("M84").Value = -532789.6699999999
$ws.Range("N84").Value = -60583

$ws.Range("H96").Value = 3345.6072
$ws.Range("I96").Value = 1461.375
$ws.Range("J96").Value = 4099.3
$ws.Range("K96").Value = 1461.375
$ws.Range("L96").Value = 4099.3
$ws.Range("M96").Value = -88.375
$ws.Range("N96").Value = -6845.3

$ws.Range("H136").Value = 35716576
$ws.Range("I136").Value = 45457000
$ws.Range("K136").Value = 136371000
$ws.Range("M136").Value = -136368450
